$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2021-09-23"

# Update the header cell text for the current-month column
$ws.Range("B1").Value = "September 2021 (through September 23)"

# "Little Village" overtook "Kenwood" in incident count and the two swap
# table rows (same row positions, neighborhood labels exchanged)
$ws.Range("A8").Value = "Little Village"
$ws.Range("A9").Value = "Kenwood"

# Per-cell count updates (includes the swapped row 8/9 data plus the new
# 2021-10-01 incident increments elsewhere on the sheet)
$ws.Range("AL2").Value = 5
$ws.Range("T2").Value = 4
$ws.Range("AC6").Value = 3
$ws.Range("AA8").Value = 2
$ws.Range("AD8").Value = 1
$ws.Range("AF8").Value = ""
$ws.Range("AG8").Value = 2
$ws.Range("AH8").Value = 1
$ws.Range("AI8").Value = 3
$ws.Range("AJ8").Value = ""
$ws.Range("AM8").Value = 1
$ws.Range("AP8").Value = 1
$ws.Range("AQ8").Value = 1
$ws.Range("AR8").Value = ""
$ws.Range("AS8").Value = 2
$ws.Range("AV8").Value = ""
$ws.Range("AX8").Value = 2
$ws.Range("B8").Value = 2
$ws.Range("BA8").Value = 2
$ws.Range("BB8").Value = 2
$ws.Range("BF8").Value = ""
$ws.Range("BG8").Value = 1
$ws.Range("BJ8").Value = 1
$ws.Range("D8").Value = 2
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = ""
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 4
$ws.Range("N8").Value = ""
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 3
$ws.Range("S8").Value = ""
$ws.Range("T8").Value = 2
$ws.Range("V8").Value = 1
$ws.Range("W8").Value = 2
$ws.Range("X8").Value = ""
$ws.Range("Y8").Value = ""
$ws.Range("AA9").Value = ""
$ws.Range("AD9").Value = ""
$ws.Range("AF9").Value = 1
$ws.Range("AG9").Value = ""
$ws.Range("AH9").Value = ""
$ws.Range("AI9").Value = 2
$ws.Range("AJ9").Value = 1
$ws.Range("AM9").Value = 2
$ws.Range("AP9").Value = ""
$ws.Range("AQ9").Value = ""
$ws.Range("AR9").Value = 1
$ws.Range("AS9").Value = ""
$ws.Range("AV9").Value = 2
$ws.Range("AX9").Value = ""
$ws.Range("B9").Value = 1
$ws.Range("BA9").Value = ""
$ws.Range("BB9").Value = ""
$ws.Range("BF9").Value = 1
$ws.Range("BG9").Value = ""
$ws.Range("BJ9").Value = ""
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 1
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 1
$ws.Range("J9").Value = 8
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = 4
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = ""
$ws.Range("P9").Value = ""
$ws.Range("Q9").Value = ""
$ws.Range("S9").Value = 1
$ws.Range("T9").Value = ""
$ws.Range("V9").Value = 2
$ws.Range("W9").Value = ""
$ws.Range("X9").Value = 1
$ws.Range("Y9").Value = 1
$ws.Range("BD10").Value = 1
$ws.Range("B12").Value = 2
$ws.Range("AL15").Value = 2
$ws.Range("T18").Value = 2
$ws.Range("K23").Value = 4
$ws.Range("AC29").Value = 1
$ws.Range("AL29").Value = 2
$ws.Range("K33").Value = 2
$ws.Range("AL39").Value = 2
$ws.Range("BD54").Value = 1
$ws.Range("AL55").Value = 2
$ws.Range("B65").Value = 2
$ws.Range("AC71").Value = 1
$ws.Range("AL85").Value = 1
$ws.Range("K97").Value = 1
$ws.Range("AC99").Value = 3
